$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$s.Shapes.Item(3).TextFrame.TextRange.Text = "Why this Building Block?"
